$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) stores values as literal text in this workbook
# (e.g. "578.00", "63.472.31") rather than as numbers. Excel's COM layer
# auto-converts numeric-looking strings assigned via .Value into real
# numbers, which would silently lose formatting (trailing zeros, the
# thousands-dot style used here, etc.). Force those cells to Text format
# first so the literal strings round-trip exactly, then restore General.
$textCells = @(
    "D2", "D3", "D5", "D6", "D7", "D8", "D9", "D10", "D13", "D16", "D17", "D18", "D19", "D20", "D22", "D24", "D25", "D26", "D27", "D28", "D29", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D47", "D48", "D49", "D51"
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '63.498.11'
$ws.Range('E2').Value = '  +5.84%  '

$ws.Range('D3').Value = '3.395.49'
$ws.Range('E3').Value = '  +6.51%  '

$ws.Range('E4').Value = '  +0.02%  '

$ws.Range('D5').Value = '577.01'
$ws.Range('E5').Value = '  +7.53%  '

$ws.Range('D6').Value = '155.73'
$ws.Range('E6').Value = '  +7.30%  '

$ws.Range('D7').Value = '0.999'
$ws.Range('E7').Value = '  -0.10%  '

$ws.Range('D8').Value = '3.399.89'
$ws.Range('E8').Value = '  +6.52%  '

$ws.Range('D9').Value = '0.533'
$ws.Range('E9').Value = '  +0.04%  '

$ws.Range('D10').Value = '7.51'
$ws.Range('E10').Value = '  +2.22%  '

$ws.Range('E11').Value = '  +7.13%  '

$ws.Range('E12').Value = '  +0.81%  '

$ws.Range('D13').Value = '3.981.46'
$ws.Range('E13').Value = '  +6.54%  '

$ws.Range('E14').Value = '  +0.33%  '

$ws.Range('E15').Value = '  +6.67%  '

$ws.Range('D16').Value = '27.08'
$ws.Range('E16').Value = '  +4.92%  '

$ws.Range('D17').Value = '63.603.22'
$ws.Range('E17').Value = '  +6.03%  '

$ws.Range('D18').Value = '3.400.73'
$ws.Range('E18').Value = '  +6.15%  '

$ws.Range('D19').Value = '6.37'
$ws.Range('E19').Value = '  +2.29%  '

$ws.Range('D20').Value = '14.04'
$ws.Range('E20').Value = '  +6.03%  '

$ws.Range('E21').Value = '  +3.07%  '

$ws.Range('D22').Value = '387.46'
$ws.Range('E22').Value = '  +4.95%  '

$ws.Range('E23').Value = '  -0.17%  '

$ws.Range('D24').Value = '0.535'
$ws.Range('E24').Value = '  +2.16%  '

$ws.Range('D25').Value = '70.74'
$ws.Range('E25').Value = '  +1.89%  '

$ws.Range('D26').Value = '9.52'
$ws.Range('E26').Value = '  +11.71%  '

$ws.Range('B27').Value = 'Kaspa'
$ws.Range('C27').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D27').Value = '0.180'
$ws.Range('E27').Value = '  +6.78%  '

$ws.Range('B28').Value = 'PEPE'
$ws.Range('C28').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D28').Value = '0.0000103'
$ws.Range('E28').Value = '  +18.21%  '

$ws.Range('D29').Value = '0.998'
$ws.Range('E29').Value = '  -0.42%  '

$ws.Range('E30').Value = '  +7.71%  '

$ws.Range('D31').Value = '6.51'
$ws.Range('E31').Value = '  +6.59%  '

$ws.Range('D32').Value = '1.35'
$ws.Range('E32').Value = '  +13.39%  '

$ws.Range('B33').Value = 'NEARProtocol'
$ws.Range('C33').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D33').Value = '5.64'
$ws.Range('E33').Value = '  +7.00%  '

$ws.Range('B34').Value = 'EthereumClassic'
$ws.Range('C34').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D34').Value = '23.10'
$ws.Range('E34').Value = '  +2.75%  '

$ws.Range('D35').Value = '6.71'
$ws.Range('E35').Value = '  +2.06%  '

$ws.Range('D36').Value = '1.49'
$ws.Range('E36').Value = '  +9.80%  '

$ws.Range('D37').Value = '158.24'
$ws.Range('E37').Value = '  +0.30%  '

$ws.Range('B38').Value = 'EnergySwap'
$ws.Range('C38').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D38').Value = '27.54'
$ws.Range('E38').Value = '  +4.22%  '

$ws.Range('B39').Value = 'Stacks'
$ws.Range('C39').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D39').Value = '1.87'
$ws.Range('E39').Value = '  +10.83%  '

$ws.Range('D40').Value = '0.0762'
$ws.Range('E40').Value = '  +7.66%  '

$ws.Range('D41').Value = '2.903.23'
$ws.Range('E41').Value = '  +4.08%  '

$ws.Range('D42').Value = '0.0322'
$ws.Range('E42').Value = '  +4.32%  '

$ws.Range('D43').Value = '0.763'
$ws.Range('E43').Value = '  +6.25%  '

$ws.Range('D44').Value = '41.42'
$ws.Range('E44').Value = '  +4.12%  '

$ws.Range('E45').Value = '  +2.17%  '

$ws.Range('E46').Value = '  +8.44%  '

$ws.Range('D47').Value = '3.444.27'
$ws.Range('E47').Value = '  +6.68%  '

$ws.Range('D48').Value = '22.33'
$ws.Range('E48').Value = '  +8.39%  '

$ws.Range('D49').Value = '299.59'
$ws.Range('E49').Value = '  +14.30%  '

$ws.Range('E50').Value = '  -1.54%  '

$ws.Range('D51').Value = '6.32'
$ws.Range('E51').Value = '  +2.86%  '

# Restore General number format on the cells we forced to Text,
# now that the literal text values are safely stored.
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "General"
}

Write-Output "applied edits"